# Updated Reporter logics and Updated code for Few testcases
#
# Each data sheet (Login / Create / Find) gets two new leading columns
# (S.No / a descriptive "Test Description"-ish column) in front of the
# existing Uname/pwd login columns; on Create/Find the old
# companyName/firstName/lastName lead-creation columns are retired in
# favor of the single description column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Login"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Login")

# shift the existing Uname/pwd header+data from A:B into C:D (same style)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C2:D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C1").Value = "Uname"
$ws.Range("D1").Value = "pwd"
$ws.Range("C2").Value = "DemoSalesManager"
$ws.Range("D2").Value = "crmsfa"
$ws.Range("C3").Value = "DemoCSR"
$ws.Range("D3").Value = "crmsfa"

$ws.Range("A1").Value = "S.No"
$ws.Range("A2").Value = "TC_0001"
$ws.Range("A3").Value = "TC_0002"
$ws.Range("B1").Value = "Test Description"
$ws.Range("B2").Value = "Verify User is able to login"
$ws.Range("B3").Value = "Verify User is able to login and logout"

$ws.Columns.Item(1).ColumnWidth = 7.3
$ws.Columns.Item(2).ColumnWidth = 34
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 6

$ws.Range("C7").Select()

# ---------------------------------------------------------------------
# Sheet "Create"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Create")

# shift the existing Uname/pwd header+data from A:B into C:D (same style)
$ws.Range("A1:B3").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# drop the retired companyName/firstName/lastName column
$ws.Range("E1:E3").Clear()

$ws.Range("A1").Value = "S.No"
$ws.Range("A2").Value = "TC_0001"
$ws.Range("A3").Value = "TC_0002"
$ws.Range("B1").Value = "Uname"
$ws.Range("B2").Value = "Verify User is able to create lead using DemoSalesManager login"
$ws.Range("B3").Value = "Verify User is able to create lead using DemoCSR login"

$ws.Columns.Item(1).ColumnWidth = 7.3
$ws.Columns.Item(2).ColumnWidth = 58.5
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 6

$ws.Range("B13").Select()

# ---------------------------------------------------------------------
# Sheet "Find"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Find")

# shift the existing Uname/pwd header+data from A:B into C:D (same style)
$ws.Range("A1:B3").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# drop the retired companyName/firstName/lastName column
$ws.Range("E1:E3").Clear()

$ws.Range("A1").Value = "S.No"
$ws.Range("A2").Value = "TC_0001"
$ws.Range("A3").Value = "TC_0002"
$ws.Range("B1").Value = "Uname"
$ws.Range("B2").Value = "Verify User is able to create and find lead using DemoSalesManager login"
$ws.Range("B3").Value = "Verify User is able to create and find lead using DemoCSR login"

$ws.Columns.Item(2).ColumnWidth = 66.5
$ws.Columns.Item(3).ColumnWidth = 18

$ws.Range("B7").Select()

# "Find" remains the active/visible sheet, matching tabSelected in the diff
$ws.Activate()
